$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Value = 2090607
$ws.Range("C4").Value = 906
$ws.Range("D4").Value = 816463
$ws.Range("E4").Value = 1158077
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 33
$ws.Range("H4").Value = 116067

$ws.Range("B7").Value = 300821
$ws.Range("C7").Value = 2538
$ws.Range("D7").Value = 149035
$ws.Range("E7").Value = 143255
$ws.Range("F7").Value = 0
$ws.Range("G7").Value = 30
$ws.Range("H7").Value = 8531

$ws.Range("B8").Value = 292950
$ws.Range("C8").Value = 1541
$ws.Range("D8").Value = 0
$ws.Range("E8").Value = 0
$ws.Range("F8").Value = 0
$ws.Range("G8").Value = 202
$ws.Range("H8").Value = 41481

$ws.Range("B12").Value = 186920
$ws.Range("C12").Value = 125
$ws.Range("D12").Value = 171600
$ws.Range("E12").Value = 6467
$ws.Range("F12").Value = 0
$ws.Range("G12").Value = 2
$ws.Range("H12").Value = 8853

$ws.Range("B19").Value = 119942
$ws.Range("C19").Value = 3921
$ws.Range("D19").Value = 81029
$ws.Range("E19").Value = 38020
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 36
$ws.Range("H19").Value = 893

$ws.Range("A27").Value = 'Suecia'
$ws.Range("B27").Value = 49684
$ws.Range("C27").Value = 205
$ws.Range("D27").Value = 0
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 0
$ws.Range("G27").Value = 40
$ws.Range("H27").Value = 4854

$ws.Range("A28").Value = 'Paises Bajos'
$ws.Range("B28").Value = 48461
$ws.Range("C28").Value = 210
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("F28").Value = 0
$ws.Range("G28").Value = 9
$ws.Range("H28").Value = 6053

$ws.Range("B35").Value = 36180
$ws.Range("C35").Value = 270
$ws.Range("D35").Value = 22200
$ws.Range("E35").Value = 12475
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 1
$ws.Range("H35").Value = 1505

$ws.Range("B51").Value = 17064
$ws.Range("C51").Value = 30
$ws.Range("D51").Value = 15985
$ws.Range("E51").Value = 404
$ws.Range("F51").Value = 0
$ws.Range("G51").Value = 1
$ws.Range("H51").Value = 675

$ws.Range("B57").Value = 12175
$ws.Range("C57").Value = 73
$ws.Range("D57").Value = 11348
$ws.Range("E57").Value = 575
$ws.Range("F57").Value = 0
$ws.Range("G57").Value = 0
$ws.Range("H57").Value = 252

$ws.Range("B58").Value = 12099
$ws.Range("C58").Value = 64
$ws.Range("D58").Value = 10993
$ws.Range("E58").Value = 512
$ws.Range("F58").Value = 0
$ws.Range("G58").Value = 1
$ws.Range("H58").Value = 594

$ws.Range("B64").Value = 9218
$ws.Range("C64").Value = 336
$ws.Range("D64").Value = 5116
$ws.Range("E64").Value = 3989
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 5
$ws.Range("H64").Value = 113

$ws.Range("B66").Value = 8613
$ws.Range("C66").Value = 5
$ws.Range("D66").Value = 8138
$ws.Range("E66").Value = 233
$ws.Range("F66").Value = 0
$ws.Range("G66").Value = 0
$ws.Range("H66").Value = 242

$ws.Range("A75").Value = 'Tayikistan'
$ws.Range("B75").Value = 4902
$ws.Range("C75").Value = 68
$ws.Range("D75").Value = 3158
$ws.Range("E75").Value = 1695
$ws.Range("F75").Value = 0
$ws.Range("G75").Value = 0
$ws.Range("H75").Value = 49

$ws.Range("A76").Value = 'Senegal'
$ws.Range("B76").Value = 4851
$ws.Range("C76").Value = 92
$ws.Range("D76").Value = 3100
$ws.Range("E76").Value = 1695
$ws.Range("F76").Value = 0
$ws.Range("G76").Value = 1
$ws.Range("H76").Value = 56

$ws.Range("A77").Value = 'Uzbekistan'
$ws.Range("B77").Value = 4837
$ws.Range("C77").Value = 96
$ws.Range("D77").Value = 3700
$ws.Range("E77").Value = 1118
$ws.Range("F77").Value = 0
$ws.Range("G77").Value = 0
$ws.Range("H77").Value = 19

$ws.Range("B78").Value = 4637
$ws.Range("C78").Value = 122
$ws.Range("D78").Value = 580
$ws.Range("E78").Value = 3956
$ws.Range("F78").Value = 0
$ws.Range("G78").Value = 3
$ws.Range("H78").Value = 101

$ws.Range("B88").Value = 3305
$ws.Range("C88").Value = 90
$ws.Range("D88").Value = 1164
$ws.Range("E88").Value = 2045
$ws.Range("F88").Value = 0
$ws.Range("G88").Value = 4
$ws.Range("H88").Value = 96

$ws.Range("A92").Value = 'Etiopia'
$ws.Range("B92").Value = 2915
$ws.Range("C92").Value = 245
$ws.Range("D92").Value = 451
$ws.Range("E92").Value = 2417
$ws.Range("F92").Value = 0
$ws.Range("G92").Value = 7
$ws.Range("H92").Value = 47

$ws.Range("A93").Value = 'Bosnia y Herzegovina'
$ws.Range("B93").Value = 2893
$ws.Range("C93").Value = 61
$ws.Range("D93").Value = 2119
$ws.Range("E93").Value = 611
$ws.Range("F93").Value = 0
$ws.Range("G93").Value = 2
$ws.Range("H93").Value = 163

$ws.Range("A94").Value = 'Venezuela'
$ws.Range("B94").Value = 2814
$ws.Range("C94").Value = 0
$ws.Range("D94").Value = 487
$ws.Range("E94").Value = 2304
$ws.Range("F94").Value = 0
$ws.Range("G94").Value = 0
$ws.Range("H94").Value = 23

$ws.Range("B106").Value = 1752
$ws.Range("C106").Value = 30
$ws.Range("D106").Value = 1023
$ws.Range("E106").Value = 628
$ws.Range("F106").Value = 0
$ws.Range("G106").Value = 4
$ws.Range("H106").Value = 101

$ws.Range("A122").Value = 'Sierra Leona'
$ws.Range("B122").Value = 1103
$ws.Range("C122").Value = 18
$ws.Range("D122").Value = 648
$ws.Range("E122").Value = 404
$ws.Range("F122").Value = 0
$ws.Range("G122").Value = 1
$ws.Range("H122").Value = 51

$ws.Range("A123").Value = 'Letonia'
$ws.Range("B123").Value = 1096
$ws.Range("C123").Value = 2
$ws.Range("D123").Value = 818
$ws.Range("E123").Value = 251
$ws.Range("F123").Value = 0
$ws.Range("G123").Value = 1
$ws.Range("H123").Value = 27

$ws.Range("A124").Value = 'Tunez'
$ws.Range("B124").Value = 1087
$ws.Range("C124").Value = 0
$ws.Range("D124").Value = 989
$ws.Range("E124").Value = 49
$ws.Range("F124").Value = 0
$ws.Range("G124").Value = 0
$ws.Range("H124").Value = 49

$ws.Range("B128").Value = 892
$ws.Range("C128").Value = 0
$ws.Range("D128").Value = 791
$ws.Range("E128").Value = 48
$ws.Range("F128").Value = 0
$ws.Range("G128").Value = 0
$ws.Range("H128").Value = 53

$ws.Range("A154").Value = 'Benin'
$ws.Range("B154").Value = 388
$ws.Range("C154").Value = 83
$ws.Range("D154").Value = 217
$ws.Range("E154").Value = 166
$ws.Range("F154").Value = 0
$ws.Range("G154").Value = 1
$ws.Range("H154").Value = 5

$ws.Range("A155").Value = 'Mauricio'
$ws.Range("B155").Value = 337
$ws.Range("C155").Value = 0
$ws.Range("D155").Value = 325
$ws.Range("E155").Value = 2
$ws.Range("F155").Value = 0
$ws.Range("G155").Value = 0
$ws.Range("H155").Value = 10

$ws.Range("A156").Value = 'Isla de Man'
$ws.Range("B156").Value = 336
$ws.Range("C156").Value = 0
$ws.Range("D156").Value = 312
$ws.Range("E156").Value = 0
$ws.Range("F156").Value = 0
$ws.Range("G156").Value = 0
$ws.Range("H156").Value = 24

$ws.Range("B157").Value = 333
$ws.Range("C157").Value = 1
$ws.Range("D157").Value = 323
$ws.Range("E157").Value = 10
$ws.Range("F157").Value = 0
$ws.Range("G157").Value = 0
$ws.Range("H157").Value = 0

$ws.Range("A158").Value = 'Zimbabue'
$ws.Range("B158").Value = 332
$ws.Range("C158").Value = 0
$ws.Range("D158").Value = 51
$ws.Range("E158").Value = 277
$ws.Range("F158").Value = 0
$ws.Range("G158").Value = 0
$ws.Range("H158").Value = 4

$ws.Range("A159").Value = 'Montenegro'
$ws.Range("B159").Value = 324
$ws.Range("C159").Value = 0
$ws.Range("D159").Value = 315
$ws.Range("E159").Value = 0
$ws.Range("F159").Value = 0
$ws.Range("G159").Value = 0
$ws.Range("H159").Value = 4
